$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the visible top-left of the sheet view up a couple of rows
# (was B13, now B11).
try {
    $excel.ActiveWindow.ScrollRow = 11
    $excel.ActiveWindow.ScrollColumn = 2
} catch {
}

# Row 22: log the hours worked and the activity summary for the final
# entry (session-closing bug fixes / test updates).
$ws.Range("E22").Value = 6
$ws.Range("F22").Value = "Finalized sesssion closing. Bug fixes. Cleaned up unused code. Added some documentation. Migrated to a different Stored Procedure to close sessions, changed tests accordingly."

# Keep the active cell/selection on F22, matching the saved view state.
$ws.Range("F22").Select()
